# rCGH.pptx update: drop the duplicate "Rescaling" slide and nudge the
# pipeline diagram's "re-echelonnees" box to its new spot/size.

$p = $ppt.ActivePresentation

# --- 1) Delete the redundant "Rescaling" slide (position 5: the simple
#        adjustSignal()/Agilent-only slide that duplicated slide 4). ---
$p.Slides.Item(5).Delete()

# --- 2) Resize/reposition the "Rectangle 15" box on the Pipeline slide
#        (slide 2) that represents the re-scaled-data step. ---
$s2 = $p.Slides.Item(2)
$rect = $s2.Shapes.Item("Rectangle 15")
$rect.Left = 4172545 / 12700
$rect.Top = 4047262 / 12700
$rect.Width = 1019691 / 12700
$rect.Height = 813496 / 12700
